# Rename the diff-report column headers so they carry the respective
# input-file "format version" as a suffix instead of the generic
# "_old" / "_new" markers, wrap the header row (and all data) in a
# native Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Column layout (A1:U1):
#   A-J  : "<Label>_old"  -> "<Label>_FV2410"
#   K    : "diff"         -> unchanged
#   L-U  : "<Label>_new"  -> "<Label>_FV2504"
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into a native Table (Table1) ------------------
$tableRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"
# No explicit table style (matches the plain, unstyled Table1 definition).
$tbl.TableStyle = ""

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
